# Update countries & provincias Spain
# - Honduras overtakes Bielorrusia in ranking (its case counts grew),
#   so the two countries swap places (row 53 / row 54).
# - Several other countries (rows 26,30,37,44,71,153,197) get refreshed
#   case-count figures as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 now holds Honduras (previously Bielorrusia), with updated figures.
$ws.Range("A53").Value = "Honduras"
$ws.Range("B53").Value = 87594
$ws.Range("C53").Value = 903
$ws.Range("D53").Value = 34662
$ws.Range("E53").Value = 50369
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 7
$ws.Range("H53").Value = 2563

# Row 54 now holds Bielorrusia (previously Honduras), keeping its old figures.
$ws.Range("A54").Value = "Bielorrusia"
$ws.Range("B54").Value = 87063
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 79429
$ws.Range("E54").Value = 6709
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 925

# Pakistan (row 26)
$ws.Range("B26").Value = 323019
$ws.Range("C26").Value = 567
$ws.Range("D26").Value = 307069
$ws.Range("E26").Value = 9296
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 6654

# Belgica (row 30)
$ws.Range("B30").Value = 213115
$ws.Range("C30").Value = 10964
$ws.Range("D30").Value = 21074
$ws.Range("E30").Value = 181649
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 33
$ws.Range("H30").Value = 10392

# Bolivia (row 37)
$ws.Range("B37").Value = 139710
$ws.Range("C37").Value = 148
$ws.Range("D37").Value = 104202
$ws.Range("E37").Value = 27045
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 24
$ws.Range("H37").Value = 8463

# Kazajistan (row 44)
$ws.Range("B44").Value = 109406
$ws.Range("C44").Value = 104
$ws.Range("D44").Value = 105001
$ws.Range("E44").Value = 2637

# Ghana (row 71)
$ws.Range("B71").Value = 47232
$ws.Range("C71").Value = 33
$ws.Range("D71").Value = 46578
$ws.Range("E71").Value = 344

# Belice (row 153)
$ws.Range("B153").Value = 2775
$ws.Range("C153").Value = 47
$ws.Range("D153").Value = 1648
$ws.Range("E153").Value = 1084

# Antigua y Barbuda (row 197)
$ws.Range("B197").Value = 119
$ws.Range("C197").Value = 6
$ws.Range("D197").Value = 101
$ws.Range("E197").Value = 15
